{"js": "// The template contained a Word field (fldChar begin/instrText.../fldChar end)\n// spelling out the M2Doc token \" m:null.check() \". The parser was switched to\n// TokenIteratorFieldRewriterSplit, which expects the token written out as\n// plain literal text \"{m:null.check()}\" instead of a real Word field.\n//\n// Find the paragraph that holds that field and replace its contents with\n// plain text runs spelling the same token, keeping the highlighted \"null\"\n// run (with its color) and the \"_GoBack\" bookmark in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  const fields = p.getRange().fields;\n  fields.load(\"items/code\");\n  await context.sync();\n  if (fields.items.length > 0 && fields.items[0].code.indexOf(\"null.check()\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the field paragraph containing 'null.check()'\");\n}\n\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>null</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t>.</w:t></w:r>' +\n  '<w:r><w:t>check()</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.getRange().insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The template contained a Word field (fldChar begin/instrText.../fldChar end)\n# spelling out the M2Doc token \" m:null.check() \". The parser was switched to\n# TokenIteratorFieldRewriterSplit, which expects the token written out as\n# plain literal text \"{m:null.check()}\" instead of a real Word field.\n#\n# Find the paragraph that holds that field and replace its contents with\n# plain text runs spelling the same token, keeping the highlighted \"null\"\n# run (with its color) and the \"_GoBack\" bookmark in place.\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    $fields = $candidate.Range.Fields\n    if ($fields.Count -gt 0 -and $fields.Item(1).Code.Text -like \"*null.check()*\") {\n        $target = $candidate\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the field paragraph containing 'null.check()'\"\n}\n\n$replacementXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>null</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n    '<w:r><w:t>check()</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>'\n\n$target.Range.InsertXML($replacementXml)\n"}
